{"js": "// Append a new run containing \".1\" right after the existing \"PRUEBA 1\"\n// run, in the same paragraph, with the identical run formatting\n// (sz=240 / szCs=240 / lang=es-ES). We insert it as its own <w:r> (not\n// merged text on the existing run) via a minimal Flat-OPC OOXML payload\n// so the new run keeps its own (attribute-less) <w:r> element, matching\n// the authored diff exactly.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[0];\nconst insertionPoint = paragraph.getRange(\"End\");\n\nconst flatOpcRunXml = `<?xml version=\"1.0\" standalone=\"yes\"?>\n<?mso-application progid=\"Word.Document\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr>\n                <w:sz w:val=\"240\"/>\n                <w:szCs w:val=\"240\"/>\n                <w:lang w:val=\"es-ES\"/>\n              </w:rPr>\n              <w:t>.1</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ninsertionPoint.insertOoxml(flatOpcRunXml, \"End\");\nawait context.sync();\n", "ps1": "# Append a new run containing \".1\" right after the existing \"PRUEBA 1\"\n# run, in the same paragraph, using the identical run formatting\n# (sz=240 / szCs=240 / lang=es-ES). InsertXML (fed a minimal Flat-OPC\n# WordprocessingML payload) inserts it as its own <w:r> element rather\n# than merging the text into the pre-existing run, matching the\n# authored diff exactly.\n$d = $word.ActiveDocument\n$p = $d.Paragraphs(1)\n$r = $p.Range\n$r.Collapse(0)\n\n$xml = @'\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<?mso-application progid=\"Word.Document\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr>\n                <w:sz w:val=\"240\"/>\n                <w:szCs w:val=\"240\"/>\n                <w:lang w:val=\"es-ES\"/>\n              </w:rPr>\n              <w:t>.1</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$r.InsertXML($xml, \"End\")\n"}
